# Update countries & provincias Spain
# Applies the 12-May-2020 data refresh (11:35 -> 12:05) to the "Pais" sheet:
#  - timestamp banner in A1
#  - Estados Unidos (row 4) totals
#  - Rumania (row 39) totals
#  - Etiopia / Cabo Verde (rows 140-141) swap order (Etiopia now has more
#    cases than Cabo Verde) and Etiopia's totals are refreshed

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Timestamp banner
$ws.Range("A1").Value = "Datos actualizados a 12 de Mayo de 2020 a las 12:05"

# Estados Unidos (row 4)
$ws.Range("B4").Value = 1385850
$ws.Range("C4").Value = 16
$ws.Range("E4").Value = 1041830

# Rumania (row 39)
$ws.Range("B39").Value = 15778
$ws.Range("C39").Value = 190
$ws.Range("D39").Value = 7685
$ws.Range("E39").Value = 7102
$ws.Range("F39").Value = 238

# Etiopia moves above Cabo Verde (row 140 becomes Etiopia with refreshed
# totals, row 141 becomes Cabo Verde keeping its previous totals)
$ws.Range("A140").Value = "Etiopia"
$ws.Range("B140").Value = 261
$ws.Range("C140").Value = 11
$ws.Range("D140").Value = 106
$ws.Range("E140").Value = 150
$ws.Range("H140").Value = 5

$ws.Range("A141").Value = "Cabo Verde"
$ws.Range("B141").Value = 260
$ws.Range("D141").Value = 58
$ws.Range("E141").Value = 200
$ws.Range("F141").Value = 0
$ws.Range("H141").Value = 2
